$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.011.27"
$ws.Range("E2").Value = "  +1.71%  "
$ws.Range("D3").Value = "2.644.18"
$ws.Range("E3").Value = "  +2.20%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.19"
$ws.Range("E5").Value = "  +2.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.09"
$ws.Range("E6").Value = "  +4.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.569"
$ws.Range("E8").Value = "  +1.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.61"
$ws.Range("E9").Value = "  +2.78%  "
$ws.Range("E10").Value = "  +2.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.338"
$ws.Range("E11").Value = "  +2.03%  "
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("D13").Value = "3.106.81"
$ws.Range("E13").Value = "  +1.92%  "
$ws.Range("D14").Value = "59.924.61"
$ws.Range("E14").Value = "  +1.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.98"
$ws.Range("E15").Value = "  +2.31%  "
$ws.Range("D16").Value = "2.642.96"
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000135"
$ws.Range("E17").Value = "  +1.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "344.51"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.40"
$ws.Range("E19").Value = "  +2.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.23"
$ws.Range("E20").Value = "  +1.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.43"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.45"
$ws.Range("E23").Value = "  +1.58%  "
$ws.Range("E24").Value = "  +1.72%  "
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.996"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("E27").Value = "  +3.02%  "
$ws.Range("D28").Value = "0.0₃0754"
$ws.Range("E28").Value = "  +4.62%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  +3.84%  "
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.95"
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "151.00"
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.03"
$ws.Range("E34").Value = "  +1.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.14"
$ws.Range("E35").Value = "  +1.81%  "
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("E37").Value = "  +4.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.841"
$ws.Range("E38").Value = "  +2.06%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "289.38"
$ws.Range("E39").Value = "  +7.27%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.58"
$ws.Range("E40").Value = "  +1.71%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("E42").Value = "  +0.96%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0953"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0532"
$ws.Range("E45").Value = "  +3.55%  "
$ws.Range("D46").Value = "1.978.91"
$ws.Range("E46").Value = "  +0.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.73"
$ws.Range("E47").Value = "  +3.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0225"
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("E49").Value = "  +1.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.68"
$ws.Range("E50").Value = "  -2.56%  "
$ws.Range("E51").Value = "  +0.26%  "
